$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 47917.332
$ws.Range("J93").Value = 47917.332
$ws.Range("L93").Value = 47917.332
$ws.Range("N93").Value = -52909.332
$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492
$ws.Range("H98").Value = 33862.082
$ws.Range("J98").Value = 72788.63
$ws.Range("L98").Value = 72788.63
$ws.Range("N98").Value = -75784.63
$ws.Range("H111").Value = 1899.5
$ws.Range("I111").Value = 1899.5
$ws.Range("K111").Value = 5698.5
$ws.Range("M111").Value = -2631.5
$ws.Range("H122").Value = 33862.082
$ws.Range("J122").Value = 72788.63
$ws.Range("L122").Value = 218365.89
$ws.Range("N122").Value = -223265.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27471
$ws.Range("J24").Value = 27471
$ws.Range("L24").Value = 27471
$ws.Range("N24").Value = -28219
$ws.Range("H100").Value = 27471
$ws.Range("J100").Value = 27471
$ws.Range("L100").Value = 27471
$ws.Range("N100").Value = -29635
$ws.Range("H101").Value = 41593.75
$ws.Range("J101").Value = 41593.75
$ws.Range("L101").Value = 41593.75
$ws.Range("N101").Value = -48083.75
$ws.Range("H103").Value = 37173.5
$ws.Range("J103").Value = 37173.5
$ws.Range("L103").Value = 37173.5
$ws.Range("N103").Value = -39517.5
$ws.Range("H104").Value = 42000
$ws.Range("J104").Value = 42000
$ws.Range("L104").Value = 42000
$ws.Range("N104").Value = -48988
$ws.Range("H105").Value = 49036.668
$ws.Range("J105").Value = 49036.668
$ws.Range("L105").Value = 49036.668
$ws.Range("N105").Value = -56024.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 19473.727
$ws.Range("J2").Value = 19473.727
$ws.Range("L2").Value = 19473.727
$ws.Range("N2").Value = -19699.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 36982.25
$ws.Range("J28").Value = 36982.25
$ws.Range("L28").Value = 36982.25
$ws.Range("N28").Value = -37472.25
$ws.Range("H99").Value = 1698.2632
$ws.Range("I99").Value = 1276.7
$ws.Range("J99").Value = 2166.6667
$ws.Range("K99").Value = 1276.7
$ws.Range("L99").Value = 2166.6667
$ws.Range("M99").Value = 221.3
$ws.Range("N99").Value = -5162.6667
$ws.Range("H106").Value = 44990.5
$ws.Range("J106").Value = 44990.5
$ws.Range("L106").Value = 44990.5
$ws.Range("N106").Value = -47514.5
$ws.Range("H124").Value = 44326
$ws.Range("J124").Value = 44326
$ws.Range("L124").Value = 44326
$ws.Range("N124").Value = -49236
$ws.Range("H125").Value = 34918.4
$ws.Range("J125").Value = 34918.4
$ws.Range("L125").Value = 34918.4
$ws.Range("N125").Value = -39838.4
$ws.Range("H126").Value = 1698.2632
$ws.Range("I126").Value = 1276.7
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 3830.1
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -1360.1
$ws.Range("N126").Value = -11440.0001
$ws.Range("H131").Value = 38326
$ws.Range("J131").Value = 38326
$ws.Range("L131").Value = 38326
$ws.Range("N131").Value = -48406

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 46535
$ws.Range("J98").Value = 46535
$ws.Range("L98").Value = 46535
$ws.Range("N98").Value = -52525
$ws.Range("H102").Value = 2009.2222
$ws.Range("I102").Value = 1867.6666
$ws.Range("J102").Value = 2504.6667
$ws.Range("K102").Value = 1867.6666
$ws.Range("L102").Value = 2504.6667
$ws.Range("M102").Value = -245.6666
$ws.Range("N102").Value = -5748.6667
$ws.Range("H104").Value = 32915.75
$ws.Range("J104").Value = 32915.75
$ws.Range("L104").Value = 32915.75
$ws.Range("N104").Value = -39903.75
$ws.Range("H120").Value = 39317
$ws.Range("J120").Value = 39317
$ws.Range("L120").Value = 39317
$ws.Range("N120").Value = -48993
$ws.Range("H125").Value = 34748
$ws.Range("J125").Value = 34748
$ws.Range("L125").Value = 34748
$ws.Range("N125").Value = -39668
$ws.Range("H126").Value = 6738.273
$ws.Range("I126").Value = 12189.2
$ws.Range("J126").Value = 2195.8333
$ws.Range("K126").Value = 36567.60000000001
$ws.Range("L126").Value = 6587.499899999999
$ws.Range("M126").Value = -34097.60000000001
$ws.Range("N126").Value = -11527.4999
$ws.Range("H129").Value = 41666
$ws.Range("J129").Value = 41666
$ws.Range("L129").Value = 41666
$ws.Range("N129").Value = -51666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2983.0588
$ws.Range("I7").Value = 2246.5454
$ws.Range("K7").Value = 2246.5454
$ws.Range("M7").Value = -2134.5454
$ws.Range("H105").Value = 48615
$ws.Range("J105").Value = 48615
$ws.Range("L105").Value = 48615
$ws.Range("N105").Value = -55603
$ws.Range("H109").Value = 35273
$ws.Range("J109").Value = 35273
$ws.Range("L109").Value = 35273
$ws.Range("N109").Value = -38047
$ws.Range("H117").Value = 42694
$ws.Range("J117").Value = 42694
$ws.Range("L117").Value = 42694
$ws.Range("N117").Value = -51872
$ws.Range("H122").Value = 2294.353
$ws.Range("I122").Value = 2022.6666
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 6067.9998
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -3617.9998
$ws.Range("N122").Value = -12700
$ws.Range("H123").Value = 35661.332
$ws.Range("J123").Value = 35661.332
$ws.Range("L123").Value = 35661.332
$ws.Range("N123").Value = -45461.332
$ws.Range("H126").Value = 2983.0588
$ws.Range("I126").Value = 2246.5454
$ws.Range("K126").Value = 6739.6362
$ws.Range("M126").Value = -4269.6362
$ws.Range("H127").Value = 50535
$ws.Range("J127").Value = 50535
$ws.Range("L127").Value = 50535
$ws.Range("N127").Value = -60455
$ws.Range("H131").Value = 42326
$ws.Range("J131").Value = 42326
$ws.Range("L131").Value = 42326
$ws.Range("N131").Value = -52406
$ws.Range("H136").Value = 2986.9565
$ws.Range("I136").Value = 2318.8125
$ws.Range("K136").Value = 6956.4375
$ws.Range("M136").Value = -4406.4375
$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 25123
$ws.Range("J27").Value = 25123
$ws.Range("L27").Value = 25123
$ws.Range("N27").Value = -25261
$ws.Range("H92").Value = 35476.668
$ws.Range("J92").Value = 35476.668
$ws.Range("L92").Value = 35476.668
$ws.Range("N92").Value = -40468.668
$ws.Range("H93").Value = 37146.668
$ws.Range("J93").Value = 37146.668
$ws.Range("L93").Value = 37146.668
$ws.Range("N93").Value = -42138.668
$ws.Range("H94").Value = 10021.667
$ws.Range("J94").Value = 10021.667
$ws.Range("L94").Value = 10021.667
$ws.Range("N94").Value = -11823.667
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H97").Value = 33549.332
$ws.Range("J97").Value = 33549.332
$ws.Range("L97").Value = 33549.332
$ws.Range("N97").Value = -35531.332
$ws.Range("H104").Value = 47155.5
$ws.Range("J104").Value = 47155.5
$ws.Range("L104").Value = 47155.5
$ws.Range("N104").Value = -54143.5
$ws.Range("H109").Value = 32784.5
$ws.Range("J109").Value = 32784.5
$ws.Range("L109").Value = 32784.5
$ws.Range("N109").Value = -35558.5
$ws.Range("H110").Value = 49640
$ws.Range("J110").Value = 49640
$ws.Range("L110").Value = 49640
$ws.Range("N110").Value = -57820
$ws.Range("H118").Value = 28782
$ws.Range("J118").Value = 28782
$ws.Range("L118").Value = 28782
$ws.Range("N118").Value = -32096
$ws.Range("H126").Value = 1548383.2
$ws.Range("I126").Value = 1548383.2
$ws.Range("K126").Value = 4645149.6
$ws.Range("M126").Value = -4642679.6
